# Grille de correction - Lab4 Vanderlay
# "plan du client non-remplis, tout le reste est parfait"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Comment for "plan du client suivie" requirement (row 8) is shortened:
# the "aucune section pour des jobs" part is removed, only the client
# plan remark remains.
$ws.Range("D8").Value = "plan du client non-remplis"

# Points awarded for that same requirement go up from 7 to 9 (out of 10),
# since the issue is less severe than initially noted.
$ws.Range("C8").Value = 9

# Totals (row 16) and percentage (row 17) recalculate automatically via
# their existing SUM / division formulas.

# Leave the selection on the total cell, as in the saved workbook.
$ws.Range("D16").Select()
